# Generate Report for Handoff
# - Mark Priority = "ht" for the handoff rows (7,8,9,11,12,13) in the
#   zh-cn and de-de sheets.
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for those same rows.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 13)

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-09-07 04:28:59"
}

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-09-07 04:29:10"
}

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-07 04:29:10"
}
